$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2-52 hold a "Förändrad" (changed) date that is being
# bumped by one day (serial 45202 -> 45203, i.e. 2023-10-03 -> 2023-10-04).
$range = $ws.Range("C2:C52")
$range.Value = 45203
